# Updated cover sheet / devlog entries, per "Updated cover sheet, assembled cards.jar"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in a new Development Log row (row 16) ---
# Date column got a (mistyped) text date rather than a real date value
$ws.Range("C16").Value = "22/22/2022"
# Time column (fraction of a day)
$ws.Range("D16").Value = 0.58333333333333337
# Duration column, formatted as a time value (h:mm) like the rows above it
$ws.Range("E16").Value = 0.09375
$ws.Range("E16").NumberFormat = "h:mm"
# Role columns for the two students
$ws.Range("F16").Value = "Driver"
$ws.Range("G16").Value = "Observer"

# --- Clear the (incorrectly duplicated) student-id values from rows 18-20 ---
$ws.Range("H18:I20").ClearContents()

# --- Update the saved view/selection state ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H23").Select()
